$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 509.450781049
$ws.Range("C2").Value = 3.49889665453

$ws.Range("B3").Value = 509450.781049
$ws.Range("C3").Value = 104966.8996359

$ws.Range("B4").Value = 1881758.349960691
$ws.Range("C4").Value = 482359.642241833

$ws.Range("B5").Value = 37635.16699921383
$ws.Range("C5").Value = 9647.192844836662
